# Insert a new row at row 135, shifting existing rows 135:253 down to 136:254.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("135:135").Insert()

# Populate the newly inserted row 135 with the new record's data.
$ws.Cells.Item(135, 1).Value = 5
$ws.Cells.Item(135, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(135, 3).Value = "Maule"
$ws.Cells.Item(135, 4).Value = 45040
$ws.Cells.Item(135, 5).Value = 7
$ws.Cells.Item(135, 6).Value = 100112031
$ws.Cells.Item(135, 7).Value = "Poroto verde"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 150
$ws.Cells.Item(135, 11).Value = 25000
$ws.Cells.Item(135, 12).Value = 25000
$ws.Cells.Item(135, 13).Value = 25000
$ws.Cells.Item(135, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(135, 15).Value = "Región del Maule"
$ws.Cells.Item(135, 16).Value = 1000
$ws.Cells.Item(135, 17).Value = 25
$ws.Cells.Item(135, 18).Value = "Hortaliza"
